$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListSheet")

# Header row: switch several headers to lowercase (NAME/LASTNAME stay as-is)
$ws.Range("A1").Value = "id"
$ws.Range("D1").Value = "birthdate"
$ws.Range("E1").Value = "begindate"
$ws.Range("F1").Value = "enddate"
$ws.Range("G1").Value = "ukrmark"
$ws.Range("H1").Value = "germark"
$ws.Range("I1").Value = "gender"
$ws.Range("J1").Value = "dl"
$ws.Range("K1").Value = "dn"
$ws.Range("L1").Value = "filldate"
$ws.Range("M1").Value = "hd"
$ws.Range("N1").Value = "md"

# Data rows: update name / lastname values
$ws.Range("B2").Value = "KONSTANTIN"
$ws.Range("C2").Value = "KOVALENKO"
$ws.Range("B3").Value = "DARYA"
$ws.Range("C3").Value = "BLABLABLA"

# Names got longer (uppercase, longer words) so the NAME column auto-fits again
$ws.Columns.Item(2).AutoFit()

# Restore the active selection to match the saved state
$ws.Range("D6").Select()
